$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 289, shifting existing rows 289-350 down to 290-351.
$ws.Rows.Item(289).Insert()

# Populate the newly inserted row 289 with the new record (same shape as the
# surrounding rows: Mercado ID / Mercado / Region / Fecha / Codreg / Categoria ID /
# Categoria / Variedad / Calidad / Volumen / Precio minimo / Precio maximo /
# Precio promedio ponderado / Unidad de comercializacion / Origen / Precio $/Kg /
# Kg o Unidades / Clasificacion).
$ws.Cells.Item(289, 1).Value = 10
$ws.Cells.Item(289, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(289, 3).Value = "La Araucanía"
$ws.Cells.Item(289, 4).Value = 44995
$ws.Cells.Item(289, 5).Value = 9
$ws.Cells.Item(289, 6).Value = 100112039
$ws.Cells.Item(289, 7).Value = "Ciboulette"
$ws.Cells.Item(289, 8).Value = "Sin especificar"
$ws.Cells.Item(289, 9).Value = "Primera"
$ws.Cells.Item(289, 10).Value = 35
$ws.Cells.Item(289, 11).Value = 5000
$ws.Cells.Item(289, 12).Value = 5000
$ws.Cells.Item(289, 13).Value = 5000
$ws.Cells.Item(289, 14).Value = "$/docena de atados"
$ws.Cells.Item(289, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(289, 16).Value = 1667
$ws.Cells.Item(289, 17).Value = 3
$ws.Cells.Item(289, 18).Value = "Hortaliza"
